# "Final fixes to IMperial and santiago"
#
# On the Santiago de Compostela sheet, Yolanda Prezado's record listed two
# affiliations: Inst-Curie (J2/K2) and Uni-PSL (L2/M2). The first
# affiliation (Inst-Curie) is removed, the remaining affiliation (Uni-PSL)
# is shifted up into the Affiliation code / Affiliation address columns
# (J2/K2), the now-empty second affiliation columns (L2/M2) are cleared,
# and the "Number of affiliations" count (I2) is updated from 2 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the second affiliation (Uni-PSL) into the first affiliation slot,
# replacing the removed Inst-Curie entry.
$ws.Range("J2").Value = $ws.Range("L2").Value()
$ws.Range("K2").Value = $ws.Range("M2").Value()

# Clear out the now-vacated second affiliation columns.
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()

# Only one affiliation remains for this row.
$ws.Range("I2").Value = 1

# Update the sheet's active selection to reflect where editing left off.
$ws.Range("O2").Select()
